$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Free up the "-----------" shared text now, so it is re-created fresh
# later (after the brand-new separator/row text below), matching the
# author's final shared-string order.
$ws.Range("B14").Value = ""

# --- Row 6: the "Name:Rapahel Epstein" label moves from B6 to E6 ---
$ws.Range("B6").Value = ""
$ws.Range("E6").Value = "Name:Rapahel Epstein"

# --- Table header moves up from row 9 to row 8 (keeps its special style) ---
# Grab C9's current (pre-overwrite) formatting - the "header" look - and
# stamp it onto C8 before filling in the values, so C8 ends up styled
# exactly like the old C9 was.
$ws.Range("C9").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C8").Value = "LineNumber"
$ws.Range("D8").Value = "PartNumber"
$ws.Range("E8").Value = "Description"
$ws.Range("F8").Value = "Item Type"
$ws.Range("G8").Value = "Price"

# --- Row 9: first data row (was row 10) - drop the old header styling first ---
$ws.Range("C9").Clear()
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = "ABC"
$ws.Range("E9").Value = "Very Good"
$ws.Range("F9").Value = "Hardware"
$ws.Range("G9").Value = 200.2

# --- Row 10: second data row (was row 11) ---
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = "DEF"
$ws.Range("E10").Value = "Not so good"
$ws.Range("F10").Value = "Software"
$ws.Range("G10").Value = 100.1

# --- Row 11: third data row (was row 12) ---
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = "GHI"
$ws.Range("E11").Value = "Really Good, Give Him a chance"
$ws.Range("F11").Value = "Software"
$ws.Range("G11").Value = 123.34

# --- Row 12: new separator row, only column D populated ---
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = "------"
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = ""

# --- Row 13: new data row ---
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = "JKL"
$ws.Range("E13").Value = "Seriously Man? "
$ws.Range("F13").Value = "No-ware"
$ws.Range("G13").Value = 0.01

# --- Row 14: footer separator label stays the same text/style ---
$ws.Range("B14").Value = "-----------"

# --- Row 15: new data row ---
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = "MNO"
$ws.Range("E15").Value = "This is okay okay"
$ws.Range("F15").Value = "Hardware"
$ws.Range("G15").Value = 121.12

# --- Move the active selection, matching the author's final cursor position ---
$ws.Range("D20").Select()
